$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Footer 1 (COM) -> word/footer2.xml, Pearson logo, id="2": image1.png -> image2.png
$ftr1 = $sec.Footers(1)
$shpA = $ftr1.Range.InlineShapes.Item(1)
$shpA2 = $shpA.Range.InlineShapes.Item(1)
$shpA2.Name = "image2.png"

# --- Footer 2 (COM) -> word/footer1.xml, Pearson logo, id="3": image1.png -> image2.png
$ftr2 = $sec.Footers(2)
$shpB = $ftr2.Range.InlineShapes.Item(1)
$shpB2 = $shpB.Range.InlineShapes.Item(1)
$shpB2.Name = "image2.png"

# --- Header 2 (COM) -> word/header1.xml, BTec logo, id="1": image2.jpg -> image1.jpg
$hdr2 = $sec.Headers(2)
$shpC = $hdr2.Range.InlineShapes.Item(1)
$shpC2 = $shpC.Range.InlineShapes.Item(1)
$shpC2.Name = "image1.jpg"

Write-Host "Footer1 shape name:" $shpA2.Name
Write-Host "Footer2 shape name:" $shpB2.Name
Write-Host "Header2 shape name:" $shpC2.Name
